# Applies the OOXML diff: 15 mis-ordered match rows are corrected by
# swapping their F:V (match detail) content between paired rows, and
# two new match rows (156, 157) are appended at the end of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Restore correct row order for 15 swapped match pairs ---
# (A:E - Indice/pais/torneio/temporada/data_partida - are untouched;
#  only F:V, the match result/odds/url columns, are exchanged.)
# rows 22 <-> 23
$ws.Range("F22").Value = 'Persik Kediri'
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 'Arema FC'
$ws.Range("I22").Value = 2
$ws.Range("J22").Value = 1.71
$ws.Range("K22").Value = '13/07/2023 22:12'
$ws.Range("L22").Value = 2.08
$ws.Range("M22").Value = '15/07/2023 09:52'
$ws.Range("N22").Value = 3.61
$ws.Range("O22").Value = '13/07/2023 22:12'
$ws.Range("P22").Value = 3.29
$ws.Range("Q22").Value = '15/07/2023 09:52'
$ws.Range("R22").Value = 4.32
$ws.Range("S22").Value = '13/07/2023 22:12'
$ws.Range("T22").Value = 3.66
$ws.Range("U22").Value = '15/07/2023 09:52'
$ws.Range("V22").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/persik-kediri-arema-fc/vXWQdmTf/'
$ws.Range("F23").Value = 'RANS Nusantara'
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 'Persita'
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 2.52
$ws.Range("K23").Value = '13/07/2023 22:12'
$ws.Range("L23").Value = 2.85
$ws.Range("M23").Value = '15/07/2023 09:59'
$ws.Range("N23").Value = 3.33
$ws.Range("O23").Value = '13/07/2023 22:12'
$ws.Range("P23").Value = 3.45
$ws.Range("Q23").Value = '15/07/2023 09:50'
$ws.Range("R23").Value = 2.54
$ws.Range("S23").Value = '13/07/2023 22:12'
$ws.Range("T23").Value = 2.4
$ws.Range("U23").Value = '15/07/2023 09:59'
$ws.Range("V23").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/rans-nusantara-persita/nczZfRc7/'

# rows 24 <-> 25
$ws.Range("F24").Value = 'Bali United'
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 'Madura United'
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 1.86
$ws.Range("K24").Value = '14/07/2023 02:12'
$ws.Range("L24").Value = 1.96
$ws.Range("M24").Value = '15/07/2023 13:58'
$ws.Range("N24").Value = 3.69
$ws.Range("O24").Value = '14/07/2023 02:12'
$ws.Range("P24").Value = 3.61
$ws.Range("Q24").Value = '15/07/2023 13:58'
$ws.Range("R24").Value = 3.52
$ws.Range("S24").Value = '14/07/2023 02:12'
$ws.Range("T24").Value = 3.7
$ws.Range("U24").Value = '15/07/2023 13:58'
$ws.Range("V24").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/bali-united-madura-united/INVUe7r1/'
$ws.Range("F25").Value = 'Persis Solo'
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 'Borneo'
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 2.81
$ws.Range("K25").Value = '14/07/2023 02:12'
$ws.Range("L25").Value = 3
$ws.Range("M25").Value = '15/07/2023 13:57'
$ws.Range("N25").Value = 3.38
$ws.Range("O25").Value = '14/07/2023 02:12'
$ws.Range("P25").Value = 3.49
$ws.Range("Q25").Value = '15/07/2023 13:57'
$ws.Range("R25").Value = 2.27
$ws.Range("S25").Value = '14/07/2023 02:12'
$ws.Range("T25").Value = 2.29
$ws.Range("U25").Value = '15/07/2023 13:57'
$ws.Range("V25").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/persis-solo-borneo/d8sMcTDl/'

# rows 42 <-> 43
$ws.Range("F42").Value = 'Persis Solo'
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 'Arema FC'
$ws.Range("I42").Value = 1
$ws.Range("J42").Value = 1.65
$ws.Range("K42").Value = '28/07/2023 22:12'
$ws.Range("L42").Value = 1.61
$ws.Range("M42").Value = '30/07/2023 09:55'
$ws.Range("N42").Value = 3.8
$ws.Range("O42").Value = '28/07/2023 22:12'
$ws.Range("P42").Value = 4.15
$ws.Range("Q42").Value = '30/07/2023 09:55'
$ws.Range("R42").Value = 4.48
$ws.Range("S42").Value = '28/07/2023 22:12'
$ws.Range("T42").Value = 5.06
$ws.Range("U42").Value = '30/07/2023 09:51'
$ws.Range("V42").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/persis-solo-arema-fc/tvOgLNBC/'
$ws.Range("F43").Value = 'RANS Nusantara'
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 'PSS Sleman'
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2.32
$ws.Range("K43").Value = '28/07/2023 22:12'
$ws.Range("L43").Value = 2.79
$ws.Range("M43").Value = '30/07/2023 09:55'
$ws.Range("N43").Value = 3.26
$ws.Range("O43").Value = '28/07/2023 22:12'
$ws.Range("P43").Value = 3.5
$ws.Range("Q43").Value = '30/07/2023 09:58'
$ws.Range("R43").Value = 2.75
$ws.Range("S43").Value = '28/07/2023 22:12'
$ws.Range("T43").Value = 2.42
$ws.Range("U43").Value = '30/07/2023 09:55'
$ws.Range("V43").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/rans-nusantara-pss-sleman/pCUpNqs0/'

# rows 44 <-> 45
$ws.Range("F44").Value = 'Barito Putera'
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 'Madura United'
$ws.Range("I44").Value = 2
$ws.Range("J44").Value = 2.12
$ws.Range("K44").Value = '29/07/2023 02:12'
$ws.Range("L44").Value = 2.34
$ws.Range("M44").Value = '30/07/2023 13:51'
$ws.Range("N44").Value = 3.31
$ws.Range("O44").Value = '29/07/2023 02:12'
$ws.Range("P44").Value = 3.44
$ws.Range("Q44").Value = '30/07/2023 13:51'
$ws.Range("R44").Value = 3.13
$ws.Range("S44").Value = '29/07/2023 02:12'
$ws.Range("T44").Value = 2.94
$ws.Range("U44").Value = '30/07/2023 13:51'
$ws.Range("V44").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/ps-barito-putera-madura-united/z1NcKsRI/'
$ws.Range("F45").Value = 'Persija Jakarta'
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 'Persebaya'
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1.62
$ws.Range("K45").Value = '29/07/2023 02:12'
$ws.Range("L45").Value = 1.67
$ws.Range("M45").Value = '30/07/2023 13:52'
$ws.Range("N45").Value = 3.93
$ws.Range("O45").Value = '29/07/2023 02:12'
$ws.Range("P45").Value = 3.9
$ws.Range("Q45").Value = '30/07/2023 13:52'
$ws.Range("R45").Value = 4.54
$ws.Range("S45").Value = '29/07/2023 02:12'
$ws.Range("T45").Value = 4.97
$ws.Range("U45").Value = '30/07/2023 13:52'
$ws.Range("V45").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/persija-jakarta-persebaya/2mPkM3d6/'

# rows 49 <-> 50
$ws.Range("F49").Value = 'Dewa United'
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 'Persis Solo'
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 2.28
$ws.Range("K49").Value = '02/08/2023 22:12'
$ws.Range("L49").Value = 2.23
$ws.Range("M49").Value = '04/08/2023 09:57'
$ws.Range("N49").Value = 3.27
$ws.Range("O49").Value = '02/08/2023 22:12'
$ws.Range("P49").Value = 3.53
$ws.Range("Q49").Value = '04/08/2023 09:57'
$ws.Range("R49").Value = 2.81
$ws.Range("S49").Value = '02/08/2023 22:12'
$ws.Range("T49").Value = 3.07
$ws.Range("U49").Value = '04/08/2023 09:57'
$ws.Range("V49").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/dewa-united-persis-solo/IN7O7PbA/'
$ws.Range("F50").Value = 'Persebaya'
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 'Persikabo 1973'
$ws.Range("I50").Value = 2
$ws.Range("J50").Value = 1.51
$ws.Range("K50").Value = '02/08/2023 22:12'
$ws.Range("L50").Value = 1.61
$ws.Range("M50").Value = '04/08/2023 09:57'
$ws.Range("N50").Value = 4.07
$ws.Range("O50").Value = '02/08/2023 22:12'
$ws.Range("P50").Value = 4.32
$ws.Range("Q50").Value = '04/08/2023 09:45'
$ws.Range("R50").Value = 5.42
$ws.Range("S50").Value = '02/08/2023 22:12'
$ws.Range("T50").Value = 4.92
$ws.Range("U50").Value = '04/08/2023 09:57'
$ws.Range("V50").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/persebaya-persikabo-1973/d8JF9oUc/'

# rows 51 <-> 52
$ws.Range("F51").Value = 'PSS Sleman'
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 'Persija Jakarta'
$ws.Range("I51").Value = 3
$ws.Range("J51").Value = 3.91
$ws.Range("K51").Value = '03/08/2023 02:12'
$ws.Range("L51").Value = 3.19
$ws.Range("M51").Value = '04/08/2023 13:55'
$ws.Range("N51").Value = 3.52
$ws.Range("O51").Value = '03/08/2023 02:12'
$ws.Range("P51").Value = 2.99
$ws.Range("Q51").Value = '04/08/2023 13:50'
$ws.Range("R51").Value = 1.81
$ws.Range("S51").Value = '03/08/2023 02:12'
$ws.Range("T51").Value = 2.44
$ws.Range("U51").Value = '04/08/2023 13:55'
$ws.Range("V51").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/pss-sleman-persija-jakarta/zaIJ85q4/'
$ws.Range("F52").Value = 'Borneo'
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 'RANS Nusantara'
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = 1.39
$ws.Range("K52").Value = '03/08/2023 02:12'
$ws.Range("L52").Value = 1.52
$ws.Range("M52").Value = '04/08/2023 13:56'
$ws.Range("N52").Value = 4.98
$ws.Range("O52").Value = '03/08/2023 02:12'
$ws.Range("P52").Value = 4.68
$ws.Range("Q52").Value = '04/08/2023 13:56'
$ws.Range("R52").Value = 5.82
$ws.Range("S52").Value = '03/08/2023 02:12'
$ws.Range("T52").Value = 5.33
$ws.Range("U52").Value = '04/08/2023 13:56'
$ws.Range("V52").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/borneo-rans-nusantara/WQBS6qEG/'

# rows 60 <-> 61
$ws.Range("F60").Value = 'PSIS Semarang'
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 'Arema FC'
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 1.67
$ws.Range("K60").Value = '07/08/2023 22:12'
$ws.Range("L60").Value = 1.92
$ws.Range("M60").Value = '09/08/2023 09:59'
$ws.Range("N60").Value = 3.71
$ws.Range("O60").Value = '07/08/2023 22:12'
$ws.Range("P60").Value = 3.44
$ws.Range("Q60").Value = '09/08/2023 09:59'
$ws.Range("R60").Value = 4.45
$ws.Range("S60").Value = '07/08/2023 22:12'
$ws.Range("T60").Value = 4.07
$ws.Range("U60").Value = '09/08/2023 09:58'
$ws.Range("V60").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/psis-semarang-arema-fc/4fbbK2rG/'
$ws.Range("F61").Value = 'Barito Putera'
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 'Dewa United'
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 2.05
$ws.Range("K61").Value = '07/08/2023 22:12'
$ws.Range("L61").Value = 2.09
$ws.Range("M61").Value = '09/08/2023 09:59'
$ws.Range("N61").Value = 3.49
$ws.Range("O61").Value = '07/08/2023 22:12'
$ws.Range("P61").Value = 3.55
$ws.Range("Q61").Value = '09/08/2023 09:57'
$ws.Range("R61").Value = 3.16
$ws.Range("S61").Value = '07/08/2023 22:12'
$ws.Range("T61").Value = 3.36
$ws.Range("U61").Value = '09/08/2023 09:59'
$ws.Range("V61").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/ps-barito-putera-dewa-united/zLgfLrT9/'

# rows 82 <-> 83
$ws.Range("F82").Value = 'Persik Kediri'
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 'PSIS Semarang'
$ws.Range("I82").Value = 1
$ws.Range("J82").Value = 2.11
$ws.Range("K82").Value = '23/08/2023 22:12'
$ws.Range("L82").Value = 1.86
$ws.Range("M82").Value = '25/08/2023 09:58'
$ws.Range("N82").Value = 3.28
$ws.Range("O82").Value = '23/08/2023 22:12'
$ws.Range("P82").Value = 3.69
$ws.Range("Q82").Value = '25/08/2023 09:58'
$ws.Range("R82").Value = 3.11
$ws.Range("S82").Value = '23/08/2023 22:12'
$ws.Range("T82").Value = 4.02
$ws.Range("U82").Value = '25/08/2023 09:58'
$ws.Range("V82").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/persik-kediri-psis-semarang/SjxkuwKr/'
$ws.Range("F83").Value = 'Madura United'
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 'FC Bhayangkara'
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1.76
$ws.Range("K83").Value = '23/08/2023 22:12'
$ws.Range("L83").Value = 1.59
$ws.Range("M83").Value = '25/08/2023 09:55'
$ws.Range("N83").Value = 3.57
$ws.Range("O83").Value = '23/08/2023 22:12'
$ws.Range("P83").Value = 4
$ws.Range("Q83").Value = '25/08/2023 09:55'
$ws.Range("R83").Value = 3.95
$ws.Range("S83").Value = '23/08/2023 22:12'
$ws.Range("T83").Value = 5.59
$ws.Range("U83").Value = '25/08/2023 09:55'
$ws.Range("V83").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/madura-united-fc-bhayangkara/AFRgvcZl/'

# rows 84 <-> 85
$ws.Range("F84").Value = 'Borneo'
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 'Persita'
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = 1.85
$ws.Range("K84").Value = '24/08/2023 02:12'
$ws.Range("L84").Value = 1.64
$ws.Range("M84").Value = '25/08/2023 13:51'
$ws.Range("N84").Value = 3.73
$ws.Range("O84").Value = '24/08/2023 02:12'
$ws.Range("P84").Value = 3.9
$ws.Range("Q84").Value = '25/08/2023 13:51'
$ws.Range("R84").Value = 3.51
$ws.Range("S84").Value = '24/08/2023 02:12'
$ws.Range("T84").Value = 5.27
$ws.Range("U84").Value = '25/08/2023 13:51'
$ws.Range("V84").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/borneo-persita/l6QcwHle/'
$ws.Range("F85").Value = 'Dewa United'
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 'Persija Jakarta'
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 3.11
$ws.Range("K85").Value = '24/08/2023 02:12'
$ws.Range("L85").Value = 2.6
$ws.Range("M85").Value = '25/08/2023 13:59'
$ws.Range("N85").Value = 3.17
$ws.Range("O85").Value = '24/08/2023 02:12'
$ws.Range("P85").Value = 3.11
$ws.Range("Q85").Value = '25/08/2023 13:59'
$ws.Range("R85").Value = 2.15
$ws.Range("S85").Value = '24/08/2023 02:12'
$ws.Range("T85").Value = 2.84
$ws.Range("U85").Value = '25/08/2023 13:59'
$ws.Range("V85").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/dewa-united-persija-jakarta/IuV1xy41/'

# rows 93 <-> 94
$ws.Range("F93").Value = 'Persija Jakarta'
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 'Persib Bandung'
$ws.Range("I93").Value = 1
$ws.Range("J93").Value = 1.96
$ws.Range("K93").Value = '31/08/2023 22:12'
$ws.Range("L93").Value = 1.96
$ws.Range("M93").Value = '02/09/2023 09:57'
$ws.Range("N93").Value = 3.33
$ws.Range("O93").Value = '31/08/2023 22:12'
$ws.Range("P93").Value = 3.19
$ws.Range("Q93").Value = '02/09/2023 09:57'
$ws.Range("R93").Value = 3.44
$ws.Range("S93").Value = '31/08/2023 22:12'
$ws.Range("T93").Value = 4.27
$ws.Range("U93").Value = '02/09/2023 09:57'
$ws.Range("V93").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/persija-jakarta-persib-bandung/vTdwxS9I/'
$ws.Range("F94").Value = 'RANS Nusantara'
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 'Persik Kediri'
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 2.42
$ws.Range("K94").Value = '31/08/2023 22:12'
$ws.Range("L94").Value = 2.2
$ws.Range("M94").Value = '02/09/2023 09:51'
$ws.Range("N94").Value = 3.25
$ws.Range("O94").Value = '31/08/2023 22:12'
$ws.Range("P94").Value = 3.16
$ws.Range("Q94").Value = '02/09/2023 09:51'
$ws.Range("R94").Value = 2.7
$ws.Range("S94").Value = '31/08/2023 22:12'
$ws.Range("T94").Value = 3.5
$ws.Range("U94").Value = '02/09/2023 09:51'
$ws.Range("V94").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/rans-nusantara-persik-kediri/d47Yx8fC/'

# rows 96 <-> 97
$ws.Range("F96").Value = 'Persikabo 1973'
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 'Dewa United'
$ws.Range("I96").Value = 1
$ws.Range("J96").Value = 2.6
$ws.Range("K96").Value = '01/09/2023 22:12'
$ws.Range("L96").Value = 4.06
$ws.Range("M96").Value = '03/09/2023 09:13'
$ws.Range("N96").Value = 3.1
$ws.Range("O96").Value = '01/09/2023 22:12'
$ws.Range("P96").Value = 3.54
$ws.Range("Q96").Value = '03/09/2023 09:19'
$ws.Range("R96").Value = 2.51
$ws.Range("S96").Value = '01/09/2023 22:12'
$ws.Range("T96").Value = 1.89
$ws.Range("U96").Value = '03/09/2023 09:19'
$ws.Range("V96").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/persikabo-1973-dewa-united/lv4QvUPa/'
$ws.Range("F97").Value = 'Persebaya'
$ws.Range("G97").Value = 2
$ws.Range("H97").Value = 'Borneo'
$ws.Range("I97").Value = 1
$ws.Range("J97").Value = 2.49
$ws.Range("K97").Value = '01/09/2023 22:13'
$ws.Range("L97").Value = 2.74
$ws.Range("M97").Value = '03/09/2023 09:59'
$ws.Range("N97").Value = 3.11
$ws.Range("O97").Value = '01/09/2023 22:13'
$ws.Range("P97").Value = 3.26
$ws.Range("Q97").Value = '03/09/2023 09:58'
$ws.Range("R97").Value = 2.71
$ws.Range("S97").Value = '01/09/2023 22:13'
$ws.Range("T97").Value = 2.59
$ws.Range("U97").Value = '03/09/2023 09:59'
$ws.Range("V97").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/persebaya-borneo/j7Qef9Pn/'

# rows 98 <-> 99
$ws.Range("F98").Value = 'PSS Sleman'
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 'PSM Makassar'
$ws.Range("I98").Value = 1
$ws.Range("J98").Value = 3.26
$ws.Range("K98").Value = '02/09/2023 02:12'
$ws.Range("L98").Value = 2.98
$ws.Range("M98").Value = '03/09/2023 13:52'
$ws.Range("N98").Value = 3.14
$ws.Range("O98").Value = '02/09/2023 02:12'
$ws.Range("P98").Value = 3.13
$ws.Range("Q98").Value = '03/09/2023 13:52'
$ws.Range("R98").Value = 2.14
$ws.Range("S98").Value = '02/09/2023 02:12'
$ws.Range("T98").Value = 2.48
$ws.Range("U98").Value = '03/09/2023 13:52'
$ws.Range("V98").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/pss-sleman-psm-makassar/6ZlPaV1P/'
$ws.Range("F99").Value = 'Barito Putera'
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 'Persis Solo'
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1.79
$ws.Range("K99").Value = '02/09/2023 02:12'
$ws.Range("L99").Value = 1.81
$ws.Range("M99").Value = '03/09/2023 13:59'
$ws.Range("N99").Value = 3.65
$ws.Range("O99").Value = '02/09/2023 02:12'
$ws.Range("P99").Value = 3.73
$ws.Range("Q99").Value = '03/09/2023 13:59'
$ws.Range("R99").Value = 3.83
$ws.Range("S99").Value = '02/09/2023 02:12'
$ws.Range("T99").Value = 4.19
$ws.Range("U99").Value = '03/09/2023 13:53'
$ws.Range("V99").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/ps-barito-putera-persis-solo/8ryiekAt/'

# rows 101 <-> 102
$ws.Range("F101").Value = 'Bali United'
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 'RANS Nusantara'
$ws.Range("I101").Value = 2
$ws.Range("J101").Value = 1.56
$ws.Range("K101").Value = '14/09/2023 02:12'
$ws.Range("L101").Value = 1.75
$ws.Range("M101").Value = '15/09/2023 13:58'
$ws.Range("N101").Value = 4.45
$ws.Range("O101").Value = '14/09/2023 02:12'
$ws.Range("P101").Value = 3.65
$ws.Range("Q101").Value = '15/09/2023 13:58'
$ws.Range("R101").Value = 4.38
$ws.Range("S101").Value = '14/09/2023 02:12'
$ws.Range("T101").Value = 4.71
$ws.Range("U101").Value = '15/09/2023 13:58'
$ws.Range("V101").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/bali-united-rans-nusantara/zeS7i795/'
$ws.Range("F102").Value = 'Dewa United'
$ws.Range("G102").Value = 2
$ws.Range("H102").Value = 'FC Bhayangkara'
$ws.Range("I102").Value = 2
$ws.Range("J102").Value = 1.95
$ws.Range("K102").Value = '14/09/2023 02:12'
$ws.Range("L102").Value = 1.76
$ws.Range("M102").Value = '15/09/2023 13:59'
$ws.Range("N102").Value = 3.36
$ws.Range("O102").Value = '14/09/2023 02:12'
$ws.Range("P102").Value = 3.62
$ws.Range("Q102").Value = '15/09/2023 13:59'
$ws.Range("R102").Value = 3.41
$ws.Range("S102").Value = '14/09/2023 02:12'
$ws.Range("T102").Value = 4.68
$ws.Range("U102").Value = '15/09/2023 13:58'
$ws.Range("V102").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/dewa-united-fc-bhayangkara/d4T3hmfa/'

# rows 127 <-> 128
$ws.Range("F127").Value = 'Persikabo 1973'
$ws.Range("G127").Value = 2
$ws.Range("H127").Value = 'Persis Solo'
$ws.Range("I127").Value = 2
$ws.Range("J127").Value = 2.54
$ws.Range("K127").Value = '04/10/2023 21:12'
$ws.Range("L127").Value = 3.74
$ws.Range("M127").Value = '06/10/2023 09:56'
$ws.Range("N127").Value = 3.24
$ws.Range("O127").Value = '04/10/2023 21:12'
$ws.Range("P127").Value = 3.73
$ws.Range("Q127").Value = '06/10/2023 09:58'
$ws.Range("R127").Value = 2.52
$ws.Range("S127").Value = '04/10/2023 21:12'
$ws.Range("T127").Value = 1.92
$ws.Range("U127").Value = '06/10/2023 09:56'
$ws.Range("V127").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/persikabo-1973-persis-solo/OE3fW2x4/'
$ws.Range("F128").Value = 'RANS Nusantara'
$ws.Range("G128").Value = 2
$ws.Range("H128").Value = 'PSIS Semarang'
$ws.Range("I128").Value = 1
$ws.Range("J128").Value = 2.54
$ws.Range("K128").Value = '04/10/2023 21:12'
$ws.Range("L128").Value = 3.07
$ws.Range("M128").Value = '06/10/2023 09:52'
$ws.Range("N128").Value = 3.19
$ws.Range("O128").Value = '04/10/2023 21:12'
$ws.Range("P128").Value = 3.19
$ws.Range("Q128").Value = '06/10/2023 09:52'
$ws.Range("R128").Value = 2.54
$ws.Range("S128").Value = '04/10/2023 21:12'
$ws.Range("T128").Value = 2.39
$ws.Range("U128").Value = '06/10/2023 09:52'
$ws.Range("V128").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/rans-nusantara-psis-semarang/j15nYO7i/'

# rows 134 <-> 135
$ws.Range("F134").Value = 'FC Bhayangkara'
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 'Bali United'
$ws.Range("I134").Value = 2
$ws.Range("J134").Value = 2.71
$ws.Range("K134").Value = '07/10/2023 01:12'
$ws.Range("L134").Value = 3.01
$ws.Range("M134").Value = '08/10/2023 13:51'
$ws.Range("N134").Value = 3.33
$ws.Range("O134").Value = '07/10/2023 01:12'
$ws.Range("P134").Value = 3.52
$ws.Range("Q134").Value = '08/10/2023 12:26'
$ws.Range("R134").Value = 2.31
$ws.Range("S134").Value = '07/10/2023 01:12'
$ws.Range("T134").Value = 2.27
$ws.Range("U134").Value = '08/10/2023 13:51'
$ws.Range("V134").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/fc-bhayangkara-bali-united/neDqACL1/'
$ws.Range("F135").Value = 'Persita'
$ws.Range("G135").Value = 2
$ws.Range("H135").Value = 'Persik Kediri'
$ws.Range("I135").Value = 2
$ws.Range("J135").Value = 2.18
$ws.Range("K135").Value = '07/10/2023 01:12'
$ws.Range("L135").Value = 2.07
$ws.Range("M135").Value = '08/10/2023 13:53'
$ws.Range("N135").Value = 3.22
$ws.Range("O135").Value = '07/10/2023 01:12'
$ws.Range("P135").Value = 3.22
$ws.Range("Q135").Value = '08/10/2023 13:53'
$ws.Range("R135").Value = 3.02
$ws.Range("S135").Value = '07/10/2023 01:12'
$ws.Range("T135").Value = 3.8
$ws.Range("U135").Value = '08/10/2023 13:53'
$ws.Range("V135").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/persita-persik-kediri/QuPzCYjk/'

# --- 2) Append two new match rows (156 and 157) ---
# Copy formatting (styles/number-formats) from the last existing row
# (155) down into the two new rows before writing their values.
$ws.Range("A155:V155").Copy()
$ws.Range("A156:V157").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# row 156
$ws.Range("A156").Value = 155
$ws.Range("B156").Value = 'indonesia'
$ws.Range("C156").Value = 'liga-1'
$ws.Range("D156").Value = '2023-2024'
$ws.Range("E156").Value = 45232.54166666666
$ws.Range("F156").Value = 'Borneo'
$ws.Range("G156").Value = 3
$ws.Range("H156").Value = 'Persik Kediri'
$ws.Range("I156").Value = 0
$ws.Range("J156").Value = 1.55
$ws.Range("K156").Value = '01/11/2023 01:12'
$ws.Range("L156").Value = 1.48
$ws.Range("M156").Value = '02/11/2023 12:54'
$ws.Range("N156").Value = 3.94
$ws.Range("O156").Value = '01/11/2023 01:12'
$ws.Range("P156").Value = 4.28
$ws.Range("Q156").Value = '02/11/2023 12:59'
$ws.Range("R156").Value = 4.88
$ws.Range("S156").Value = '01/11/2023 01:12'
$ws.Range("T156").Value = 6.89
$ws.Range("U156").Value = '02/11/2023 12:59'
$ws.Range("V156").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/borneo-persik-kediri/ABg8YUJP/'

# row 157
$ws.Range("A157").Value = 156
$ws.Range("B157").Value = 'indonesia'
$ws.Range("C157").Value = 'liga-1'
$ws.Range("D157").Value = '2023-2024'
$ws.Range("E157").Value = 45232.54166666666
$ws.Range("F157").Value = 'FC Bhayangkara'
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 'PSIS Semarang'
$ws.Range("I157").Value = 1
$ws.Range("J157").Value = 3.01
$ws.Range("K157").Value = '01/11/2023 01:12'
$ws.Range("L157").Value = 3.74
$ws.Range("M157").Value = '02/11/2023 12:56'
$ws.Range("N157").Value = 3.27
$ws.Range("O157").Value = '01/11/2023 01:12'
$ws.Range("P157").Value = 3.53
$ws.Range("Q157").Value = '02/11/2023 12:59'
$ws.Range("R157").Value = 2.16
$ws.Range("S157").Value = '01/11/2023 01:12'
$ws.Range("T157").Value = 1.97
$ws.Range("U157").Value = '02/11/2023 12:58'
$ws.Range("V157").Value = 'https://www.betexplorer.com/football/indonesia/liga-1/fc-bhayangkara-psis-semarang/29YnRnCa/'

